{"js": "// Update the worksheet date and the 25 division problems to the new\n// day's values. Each \"find\" string is unique within the document, so a\n// direct search + replace (one hit each) safely re-targets every run\n// without disturbing formatting (rPr/pPr stay untouched since we only\n// rewrite the text of the existing run via Range.insertText(..., \"Replace\")).\n\nconst replacements = [\n  { find: \"2025-04-12 Saturday\", replace: \"2025-04-13 Sunday\" },\n  { find: \"50\u00f77=\", replace: \"90\u00f72=\" },\n  { find: \"10\u00f76=\", replace: \"69\u00f78=\" },\n  { find: \"11\u00f72=\", replace: \"16\u00f77=\" },\n  { find: \"36\u00f75=\", replace: \"49\u00f73=\" },\n  { find: \"87\u00f77=\", replace: \"84\u00f75=\" },\n  { find: \"63\u00f76=\", replace: \"77\u00f75=\" },\n  { find: \"77\u00f74=\", replace: \"26\u00f77=\" },\n  { find: \"74\u00f74=\", replace: \"24\u00f77=\" },\n  { find: \"73\u00f72=\", replace: \"58\u00f73=\" },\n  { find: \"30\u00f73=\", replace: \"88\u00f79=\" },\n  { find: \"45\u00f72=\", replace: \"30\u00f75=\" },\n  { find: \"91\u00f78=\", replace: \"99\u00f75=\" },\n  { find: \"93\u00f75=\", replace: \"14\u00f79=\" },\n  { find: \"58\u00f72=\", replace: \"39\u00f76=\" },\n  { find: \"72\u00f73=\", replace: \"67\u00f78=\" },\n  { find: \"65\u00f77=\", replace: \"72\u00f75=\" },\n  { find: \"39\u00f73=\", replace: \"71\u00f79=\" },\n  { find: \"51\u00f74=\", replace: \"80\u00f77=\" },\n  { find: \"88\u00f75=\", replace: \"33\u00f73=\" },\n  { find: \"66\u00f78=\", replace: \"32\u00f72=\" },\n  { find: \"13\u00f78=\", replace: \"47\u00f73=\" },\n  { find: \"96\u00f72=\", replace: \"63\u00f74=\" },\n  { find: \"95\u00f72=\", replace: \"84\u00f79=\" },\n  { find: \"62\u00f77=\", replace: \"62\u00f78=\" },\n  { find: \"18\u00f74=\", replace: \"51\u00f73=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${find}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 division problems to the new\n# day's values. Every \"Find\" string is unique in the document, so a\n# simple Find/Replace-All (one hit each) safely retargets each run's\n# text without touching any other formatting.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"2025-04-12 Saturday\"; Replace = \"2025-04-13 Sunday\" },\n    @{ Find = \"50\u00f77=\"; Replace = \"90\u00f72=\" },\n    @{ Find = \"10\u00f76=\"; Replace = \"69\u00f78=\" },\n    @{ Find = \"11\u00f72=\"; Replace = \"16\u00f77=\" },\n    @{ Find = \"36\u00f75=\"; Replace = \"49\u00f73=\" },\n    @{ Find = \"87\u00f77=\"; Replace = \"84\u00f75=\" },\n    @{ Find = \"63\u00f76=\"; Replace = \"77\u00f75=\" },\n    @{ Find = \"77\u00f74=\"; Replace = \"26\u00f77=\" },\n    @{ Find = \"74\u00f74=\"; Replace = \"24\u00f77=\" },\n    @{ Find = \"73\u00f72=\"; Replace = \"58\u00f73=\" },\n    @{ Find = \"30\u00f73=\"; Replace = \"88\u00f79=\" },\n    @{ Find = \"45\u00f72=\"; Replace = \"30\u00f75=\" },\n    @{ Find = \"91\u00f78=\"; Replace = \"99\u00f75=\" },\n    @{ Find = \"93\u00f75=\"; Replace = \"14\u00f79=\" },\n    @{ Find = \"58\u00f72=\"; Replace = \"39\u00f76=\" },\n    @{ Find = \"72\u00f73=\"; Replace = \"67\u00f78=\" },\n    @{ Find = \"65\u00f77=\"; Replace = \"72\u00f75=\" },\n    @{ Find = \"39\u00f73=\"; Replace = \"71\u00f79=\" },\n    @{ Find = \"51\u00f74=\"; Replace = \"80\u00f77=\" },\n    @{ Find = \"88\u00f75=\"; Replace = \"33\u00f73=\" },\n    @{ Find = \"66\u00f78=\"; Replace = \"32\u00f72=\" },\n    @{ Find = \"13\u00f78=\"; Replace = \"47\u00f73=\" },\n    @{ Find = \"96\u00f72=\"; Replace = \"63\u00f74=\" },\n    @{ Find = \"95\u00f72=\"; Replace = \"84\u00f79=\" },\n    @{ Find = \"62\u00f77=\"; Replace = \"62\u00f78=\" },\n    @{ Find = \"18\u00f74=\"; Replace = \"51\u00f73=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Find\n    $find.Replacement.Text = $pair.Replace\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2) | Out-Null\n}\n"}
